$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Remodel" (currently in C22, the "to implement" list) has been finished,
# so move it into the "Finished Cards" list at D27, vacating C22.
$ws.Range("C22").ClearContents()
$ws.Range("D27").Value = "Remodel"

# "Workshop" (currently in C27, the "to implement" list) has also been
# finished, so move it into the "Finished Cards" list at D26, vacating C27.
$ws.Range("C27").ClearContents()
$ws.Range("D26").Value = "Workshop"

# Update the active selection to reflect the new cursor location.
$ws.Range("C22").Select()
